# Updates the cryptocurrency price/volume snapshot in the active worksheet.
# Each updated cell currently holds a text value (prices such as "40.413.24"
# and percentages such as "  +3.72%  " are stored as text, not numbers), so
# new values are written with a leading apostrophe to force Excel to keep
# storing them as text (matching the workbook's original inlineStr typing)
# instead of auto-coercing number-looking strings into numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = '40.260.43';
    "E2" = '  +3.25%  ';
    "D3" = '2.248.64';
    "E3" = '  +0.83%  ';
    "E4" = '  +0.13%  ';
    "D5" = '295.29';
    "E5" = '  -0.44%  ';
    "D6" = '87.41';
    "E6" = '  +9.06%  ';
    "E7" = '  +1.89%  ';
    "E8" = '  +0.00%  ';
    "D9" = '0.477';
    "E9" = '  +3.94%  ';
    "D10" = '31.45';
    "E10" = '  +12.91%  ';
    "D11" = '0.0800';
    "E11" = '  +3.90%  ';
    "D12" = '47.39';
    "E12" = '  +3.41%  ';
    "E13" = '  +1.01%  ';
    "D14" = '6.49';
    "E14" = '  +6.56%  ';
    "D15" = '2.605.51';
    "E15" = '  +1.53%  ';
    "D16" = '14.33';
    "E16" = '  +1.74%  ';
    "D17" = '2.261.83';
    "E17" = '  +1.15%  ';
    "D18" = '0.742';
    "E18" = '  +3.88%  ';
    "D19" = '40.221.08';
    "E19" = '  +3.35%  ';
    "D20" = '0.0' + [string][char]0x2083 + '0894';
    "E20" = '  +4.39%  ';
    "D21" = '5.86';
    "E21" = '  +2.48%  ';
    "E22" = '  +9.86%  ';
    "D23" = '65.91';
    "E23" = '  +1.70%  ';
    "D24" = '237.51';
    "E24" = '  +5.68%  ';
    "E25" = '  -0.01%  ';
    "D26" = '2.48';
    "E26" = '  +4.39%  ';
    "D27" = '1.85';
    "E27" = '  +8.14%  ';
    "D28" = '23.16';
    "E28" = '  +4.42%  ';
    "D29" = '2.16';
    "E29" = '  -1.34%  ';
    "D30" = '9.31';
    "E30" = '  +4.75%  ';
    "D31" = '33.66';
    "E31" = '  +8.27%  ';
    "D32" = '153.41';
    "E32" = '  +3.15%  ';
    "E33" = '  +0.16%  ';
    "D34" = '4.93';
    "E34" = '  +3.58%  ';
    "D35" = '0.0722';
    "E35" = '  +5.64%  ';
    "E36" = '  +2.59%  ';
    "D37" = '16.77';
    "E37" = '  +15.94%  ';
    "D38" = '0.102';
    "E38" = '  +7.28%  ';
    "E39" = '  +2.91%  ';
    "D40" = '2.74';
    "E40" = '  +2.84%  ';
    "D41" = '1.70';
    "E41" = '  +6.26%  ';
    "D42" = '3.83';
    "E42" = '  +5.28%  ';
    "D43" = '2.025.61';
    "E43" = '  +6.33%  ';
    "D44" = '2.23';
    "E44" = '  +7.80%  ';
    "D45" = '0.0273';
    "E45" = '  +7.85%  ';
    "D46" = '10.02';
    "E46" = '  +11.47%  ';
    "D47" = '16.55';
    "E47" = '  +2.32%  ';
    "D48" = '2.61';
    "E48" = '  +3.91%  ';
    "D49" = '2.484.49';
    "E49" = '  +2.00%  ';
    "D50" = '71.98';
    "E50" = '  +6.10%  ';
    "E51" = '  +15.42%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = "'" + $updates[$cellRef]
}

Write-Output "Applied cryptos update to $($updates.Count) cells"
